# "Generate Report for handoff": the handoff transform failed for the one
# tracked file in both the zh-cn and de-de targets, so:
#   - Status flips from "Ready for handoff" to "Handoff transform failed"
#     (this shows up on the Overview roll-up sheet too, since it shares
#     the same status text per language).
#   - The just-generated handoff file/link is gone, so "Latest Handoff
#     File" is cleared (cell + hyperlink removed).
#   - "Latest Handoff Datetime" and "Latest Handback DateTime" reset to
#     the zero date.
#   - "Handoff Reason" flips from "Include" to "Ignored".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the hyperlink that lives on C2 (the "Latest Handoff File"
    # link to the generated .xlf) - this also renumbers the remaining
    # hyperlink relationship ids.
    foreach ($h in @($ws.Hyperlinks)) {
        if ($h.Range.Address() -eq '$C$2') {
            $h.Delete()
        }
    }

    # Status: "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff File is now empty - drop the whole cell.
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime resets to the zero date.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Latest Handback DateTime resets to the zero date.
    $ws.Range("G2").Value = "0001-01-01 00:00:00"

    # Handoff Reason: "Include" -> "Ignored"
    $ws.Range("H2").Value = "Ignored"
}
